## Update metadata.xlsx:
##  - rename the "groups" tab to "form_level_data"
##  - remove the "item_type" column from the form_level_data table
##  - add a new "review_required" boolean column to that table
##  - misc view-state touch ups (active tab / selection) that came along
##    with the author re-saving the workbook in Excel

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename "groups" -> "form_level_data"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("groups")
$ws.Name = "form_level_data"

# ---------------------------------------------------------------------
# 2. Rework the table on that sheet: drop "item_type", append
#    "review_required" (TRUE/FALSE) as the new last column.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()                       # keep the data, drop the table wrapper

$ws.Columns.Item(2).Delete()       # removes "item_type" (col B), shifts left

$ws.Range("D1").Value = "review_required"
$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $true

$newLo = $ws.ListObjects.Add(1, $ws.Range("A1:D3"), $null, 1)
$newLo.Name = "Table2"
$newLo.TableStyle = "TableStyleLight1"

# ---------------------------------------------------------------------
# 3. View-state touch-ups that show up in the saved file: the
#    "form_level_data" sheet becomes the active tab, and its selection
#    moved to the new review_required column (D4).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("D4").Select()
